$wb = $excel.ActiveWorkbook

# --- Sheet: Match Data ---
$ws1 = $wb.Worksheets.Item("Match Data")

$ws1.Range("A2").Value = "HVGN0BW0"
$ws1.Range("C2").Value = "KGLE38K4"

$ws1.Range("A3").Value = "MJST45X9"
$ws1.Range("B3").Value = "SWQR78Z2"
$ws1.Range("C3").Value = "TBPL91M5"
$ws1.Range("D3").Value = "LCKM33Y8"

$ws1.Range("A4").Value = "DLVW67N4"
$ws1.Range("C4").Value = "EDRX29H6"

$ws1.Range("A5").Value = "AWJF82P1"
$ws1.Range("B5").Value = "JLMN56Q3"
$ws1.Range("C5").Value = "RZQW74T7"
$ws1.Range("D5").Value = "MKHY93V0"

$ws1.Range("A6").Value = "CTBR48K5"
$ws1.Range("C6").Value = "ARLZ61F9"

# --- Sheet: Instructions ---
$ws2 = $wb.Worksheets.Item("Instructions")

$ws2.Range("A3").Value = "Format Guidelines:"
$ws2.Range("A4").Value = [char]0x2022 + " Use player passport codes (e.g., HVGN0BW0, KGLE38K4)"
$ws2.Range("A5").Value = [char]0x2022 + " Leave Team_1_Player_2 and Team_2_Player_2 empty for singles matches"
$ws2.Range("A6").Value = [char]0x2022 + " Use scores like: 11, 7, 15, 13 (games to points)"
$ws2.Range("A7").Value = [char]0x2022 + " Date format: YYYY-MM-DD (e.g., 2025-01-15)"
$ws2.Range("A8").Value = [char]0x2022 + " Gender_Override: M or F (optional, for cross-gender matches)"

$ws2.Range("A10").Value = "Examples:"
$ws2.Range("A11").Value = "Singles: HVGN0BW0 vs KGLE38K4, scores 11-7"
$ws2.Range("A12").Value = "Doubles: MJST45X9/SWQR78Z2 vs TBPL91M5/LCKM33Y8, scores 11-9"
$ws2.Range("A13").Value = ""
$ws2.Range("A14").Value = "Validation will check:"
$ws2.Range("A15").Value = [char]0x2022 + " All passport codes exist in the system"
$ws2.Range("A16").Value = [char]0x2022 + " Valid score formats"
$ws2.Range("A17").Value = [char]0x2022 + " No duplicate matches"
$ws2.Range("A18").Value = [char]0x2022 + " Proper date formatting"
